# Generate Report for Handback
# Updates the localization-status report: marks handback as complete
# ("Handed back: in sync with en-US"), refreshes the Latest Handback
# DateTime stamps, and clears the stale "handback file is not latest"
# error now that the handback is current.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status column updates (was "Ready for handoff") ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- zh-cn: Latest Handback DateTime + cleared Error Detail ---
$wsZhCn.Range("K2").Value = "2016-08-31 16:57:43"
$wsZhCn.Range("P2").Value = ""

# --- de-de: Latest Handback DateTime + cleared Error Detail ---
$wsDeDe.Range("K2").Value = "2016-08-31 16:57:50"
$wsDeDe.Range("P2").Value = ""

# --- Column width refresh (status text got longer, error detail got shorter) ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
